# Fruta / hortaliza, semanal
#
# The underlying data rows (2-20) get reshuffled: each row keeps the
# columns that are constant across the whole table (Mercado ID, Mercado,
# Región, Codreg, Categoría ID, Categoría, Unidad de comercialización,
# Kg o Unidades, Clasificación) but the "record" columns (Fecha, Variedad,
# Calidad, Volumen, Precio mínimo/máximo/promedio ponderado, Origen,
# Precio $/Kg) are permuted across rows.
#
# Build a snapshot of the record columns for every row first (so the
# permutation reads are not clobbered by earlier writes), then re-write
# each row from the snapshot of its source row according to the mapping
# derived from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-record data which gets shuffled between rows.
# D=Fecha, H=Variedad, I=Calidad, J=Volumen, K=Precio minimo,
# L=Precio maximo, M=Precio promedio ponderado, O=Origen, P=Precio $/Kg
$cols = @(4, 8, 9, 10, 11, 12, 13, 15, 16)

$firstRow = 2
$lastRow = 20

# Snapshot every row's current values before making any changes.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# new row -> source (old) row, derived from the diff.
$mapping = @{
    2  = 10
    3  = 16
    4  = 15
    5  = 13
    6  = 9
    7  = 8
    8  = 4
    9  = 2
    10 = 6
    11 = 17
    12 = 18
    13 = 20
    14 = 14
    15 = 5
    16 = 11
    17 = 12
    18 = 7
    19 = 19
    20 = 3
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $data = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value2 = $data[$c]
    }
}
